# Auto-generated edit script: update market-price derived columns (H:N)
# on the leve-profit sheets, matching a scheduled price-refresh run.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Cells.Item(33, 8).Value = 11813.556   # H33: 14865.857 -> 11813.556
$ws.Cells.Item(33, 9).Value = 20324   # I33: 20292 -> 20324
$ws.Cells.Item(33, 10).Value = 1175.5   # J33: 1300.5 -> 1175.5
$ws.Cells.Item(33, 11).Value = 20324   # K33: 20292 -> 20324
$ws.Cells.Item(33, 12).Value = 1175.5   # L33: 1300.5 -> 1175.5
$ws.Cells.Item(33, 13).Value = -20095   # M33: -20063 -> -20095
$ws.Cells.Item(33, 14).Value = -1633.5   # N33: -1758.5 -> -1633.5

# Row 63
$ws.Cells.Item(63, 8).Value = 69999   # H63: 70000 -> 69999
$ws.Cells.Item(63, 9).Value = 0   # I63: 70000 -> 0
$ws.Cells.Item(63, 10).Value = 69999   # J63: 0 -> 69999
$ws.Cells.Item(63, 11).Value = 0   # K63: 70000 -> 0
$ws.Cells.Item(63, 12).Value = 69999   # L63: 0 -> 69999
$ws.Cells.Item(63, 13).Value = ""   # M63: -69376 -> None
$ws.Cells.Item(63, 14).Value = -71247   # N63: None -> -71247

# Row 66
$ws.Cells.Item(66, 8).Value = 69999   # H66: 70000 -> 69999
$ws.Cells.Item(66, 9).Value = 0   # I66: 70000 -> 0
$ws.Cells.Item(66, 10).Value = 69999   # J66: 0 -> 69999
$ws.Cells.Item(66, 11).Value = 0   # K66: 210000 -> 0
$ws.Cells.Item(66, 12).Value = 209997   # L66: 0 -> 209997
$ws.Cells.Item(66, 13).Value = ""   # M66: -206880 -> None
$ws.Cells.Item(66, 14).Value = -216237   # N66: None -> -216237

# Row 107
$ws.Cells.Item(107, 8).Value = 839   # H107: 988.1429000000001 -> 839
$ws.Cells.Item(107, 9).Value = 349.2   # I107: 402.75 -> 349.2
$ws.Cells.Item(107, 10).Value = 1451.25   # J107: 1768.6666 -> 1451.25
$ws.Cells.Item(107, 11).Value = 349.2   # K107: 402.75 -> 349.2
$ws.Cells.Item(107, 12).Value = 1451.25   # L107: 1768.6666 -> 1451.25
$ws.Cells.Item(107, 13).Value = 1570.8   # M107: 1517.25 -> 1570.8
$ws.Cells.Item(107, 14).Value = -5291.25   # N107: -5608.6666 -> -5291.25

# Row 113
$ws.Cells.Item(113, 8).Value = 6582.8887   # H113: 6313.273 -> 6582.8887
$ws.Cells.Item(113, 9).Value = 5837.25   # I113: 5591.5 -> 5837.25
$ws.Cells.Item(113, 11).Value = 5837.25   # K113: 5591.5 -> 5837.25
$ws.Cells.Item(113, 13).Value = -2583.25   # M113: -2337.5 -> -2583.25

# Row 117
$ws.Cells.Item(117, 8).Value = 79999   # H117: 0 -> 79999
$ws.Cells.Item(117, 10).Value = 79999   # J117: 0 -> 79999
$ws.Cells.Item(117, 12).Value = 79999   # L117: 0 -> 79999
$ws.Cells.Item(117, 14).Value = -89177   # N117: None -> -89177

# Row 138
$ws.Cells.Item(138, 8).Value = 3709.1738   # H138: 3756 -> 3709.1738
$ws.Cells.Item(138, 9).Value = 3714.1365   # I138: 3763.4285 -> 3714.1365
$ws.Cells.Item(138, 11).Value = 11142.4095   # K138: 11290.2855 -> 11142.4095
$ws.Cells.Item(138, 13).Value = -6002.4095   # M138: -6150.2855 -> -6002.4095

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Cells.Item(61, 8).Value = 2472.389   # H61: 2643.0715 -> 2472.389
$ws.Cells.Item(61, 9).Value = 2566.8667   # I61: 2643.0715 -> 2566.8667
$ws.Cells.Item(61, 10).Value = 2000   # J61: 0 -> 2000
$ws.Cells.Item(61, 11).Value = 2566.8667   # K61: 2643.0715 -> 2566.8667
$ws.Cells.Item(61, 12).Value = 2000   # L61: 0 -> 2000
$ws.Cells.Item(61, 13).Value = -2354.8667   # M61: -2431.0715 -> -2354.8667
$ws.Cells.Item(61, 14).Value = -2424   # N61: None -> -2424

# Row 74
$ws.Cells.Item(74, 8).Value = 2208.5   # H74: 2024.2222 -> 2208.5
$ws.Cells.Item(74, 9).Value = 2069.6   # I74: 1635.4286 -> 2069.6
$ws.Cells.Item(74, 11).Value = 2069.6   # K74: 1635.4286 -> 2069.6
$ws.Cells.Item(74, 13).Value = -1195.6   # M74: -761.4286 -> -1195.6

# Row 77
$ws.Cells.Item(77, 8).Value = 2208.5   # H77: 2024.2222 -> 2208.5
$ws.Cells.Item(77, 9).Value = 2069.6   # I77: 1635.4286 -> 2069.6
$ws.Cells.Item(77, 11).Value = 10348   # K77: 8177.143 -> 10348
$ws.Cells.Item(77, 13).Value = -5980   # M77: -3809.143 -> -5980

# Row 122
$ws.Cells.Item(122, 8).Value = 2645.7778   # H122: 1375.1428 -> 2645.7778
$ws.Cells.Item(122, 9).Value = 2202.4   # I122: 1168.16 -> 2202.4
$ws.Cells.Item(122, 10).Value = 3200   # J122: 3100 -> 3200
$ws.Cells.Item(122, 11).Value = 6607.200000000001   # K122: 3504.48 -> 6607.200000000001
$ws.Cells.Item(122, 12).Value = 9600   # L122: 9300 -> 9600
$ws.Cells.Item(122, 13).Value = -4157.200000000001   # M122: -1054.48 -> -4157.200000000001
$ws.Cells.Item(122, 14).Value = -14500   # N122: -14200 -> -14500

# Row 136
$ws.Cells.Item(136, 8).Value = 2472.389   # H136: 2643.0715 -> 2472.389
$ws.Cells.Item(136, 9).Value = 2566.8667   # I136: 2643.0715 -> 2566.8667
$ws.Cells.Item(136, 10).Value = 2000   # J136: 0 -> 2000
$ws.Cells.Item(136, 11).Value = 7700.6001   # K136: 7929.2145 -> 7700.6001
$ws.Cells.Item(136, 12).Value = 6000   # L136: 0 -> 6000
$ws.Cells.Item(136, 13).Value = -5150.6001   # M136: -5379.2145 -> -5150.6001
$ws.Cells.Item(136, 14).Value = -11100   # N136: None -> -11100

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Cells.Item(64, 8).Value = 2680.1667   # H64: 2820.1765 -> 2680.1667
$ws.Cells.Item(64, 9).Value = 1944.6   # I64: 2105.111 -> 1944.6
$ws.Cells.Item(64, 10).Value = 3599.625   # J64: 3624.625 -> 3599.625
$ws.Cells.Item(64, 11).Value = 1944.6   # K64: 2105.111 -> 1944.6
$ws.Cells.Item(64, 12).Value = 3599.625   # L64: 3624.625 -> 3599.625
$ws.Cells.Item(64, 13).Value = -1719.6   # M64: -1880.111 -> -1719.6
$ws.Cells.Item(64, 14).Value = -4049.625   # N64: -4074.625 -> -4049.625

# Row 67
$ws.Cells.Item(67, 8).Value = 2680.1667   # H67: 2820.1765 -> 2680.1667
$ws.Cells.Item(67, 9).Value = 1944.6   # I67: 2105.111 -> 1944.6
$ws.Cells.Item(67, 10).Value = 3599.625   # J67: 3624.625 -> 3599.625
$ws.Cells.Item(67, 11).Value = 1944.6   # K67: 2105.111 -> 1944.6
$ws.Cells.Item(67, 12).Value = 3599.625   # L67: 3624.625 -> 3599.625
$ws.Cells.Item(67, 13).Value = -1164.6   # M67: -1325.111 -> -1164.6
$ws.Cells.Item(67, 14).Value = -5159.625   # N67: -5184.625 -> -5159.625

# Row 138
$ws.Cells.Item(138, 8).Value = 98920   # H138: 98898.664 -> 98920
$ws.Cells.Item(138, 10).Value = 98920   # J138: 98898.664 -> 98920
$ws.Cells.Item(138, 12).Value = 98920   # L138: 98898.664 -> 98920
$ws.Cells.Item(138, 14).Value = -109200   # N138: -109178.664 -> -109200

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 45
$ws.Cells.Item(45, 8).Value = 21750   # H45: 0 -> 21750
$ws.Cells.Item(45, 9).Value = 23500   # I45: 0 -> 23500
$ws.Cells.Item(45, 10).Value = 20000   # J45: 0 -> 20000
$ws.Cells.Item(45, 11).Value = 23500   # K45: 0 -> 23500
$ws.Cells.Item(45, 12).Value = 20000   # L45: 0 -> 20000
$ws.Cells.Item(45, 13).Value = -22907   # M45: None -> -22907
$ws.Cells.Item(45, 14).Value = -21186   # N45: None -> -21186

# Row 58
$ws.Cells.Item(58, 8).Value = 45473.824   # H58: 38766.258 -> 45473.824
$ws.Cells.Item(58, 9).Value = 49533.24   # I58: 41639.56 -> 49533.24
$ws.Cells.Item(58, 11).Value = 49533.24   # K58: 41639.56 -> 49533.24
$ws.Cells.Item(58, 13).Value = -49330.24   # M58: -41436.56 -> -49330.24

# Row 105
$ws.Cells.Item(105, 8).Value = 1044.5555   # H105: 1111.375 -> 1044.5555
$ws.Cells.Item(105, 9).Value = 1044.5555   # I105: 1111.375 -> 1044.5555
$ws.Cells.Item(105, 11).Value = 1044.5555   # K105: 1111.375 -> 1044.5555
$ws.Cells.Item(105, 13).Value = 702.4445000000001   # M105: 635.625 -> 702.4445000000001

# Row 122
$ws.Cells.Item(122, 8).Value = 2739.875   # H122: 2845.7144 -> 2739.875
$ws.Cells.Item(122, 9).Value = 2702.7144   # I122: 2820 -> 2702.7144
$ws.Cells.Item(122, 11).Value = 8108.1432   # K122: 8460 -> 8108.1432
$ws.Cells.Item(122, 13).Value = -5658.1432   # M122: -6010 -> -5658.1432

# Row 132
$ws.Cells.Item(132, 8).Value = 2901   # H132: 3330 -> 2901
$ws.Cells.Item(132, 10).Value = 2740   # J132: 4000 -> 2740
$ws.Cells.Item(132, 12).Value = 8220   # L132: 12000 -> 8220
$ws.Cells.Item(132, 14).Value = -13280   # N132: -17060 -> -13280

# Row 136
$ws.Cells.Item(136, 8).Value = 45473.824   # H136: 38766.258 -> 45473.824
$ws.Cells.Item(136, 9).Value = 49533.24   # I136: 41639.56 -> 49533.24
$ws.Cells.Item(136, 11).Value = 148599.72   # K136: 124918.68 -> 148599.72
$ws.Cells.Item(136, 13).Value = -146049.72   # M136: -122368.68 -> -146049.72

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 2755.182   # H102: 2848.7 -> 2755.182
$ws.Cells.Item(102, 9).Value = 2445.375   # I102: 2491.7334 -> 2445.375
$ws.Cells.Item(102, 10).Value = 3581.3333   # J102: 3919.6 -> 3581.3333
$ws.Cells.Item(102, 11).Value = 2445.375   # K102: 2491.7334 -> 2445.375
$ws.Cells.Item(102, 12).Value = 3581.3333   # L102: 3919.6 -> 3581.3333
$ws.Cells.Item(102, 13).Value = -823.375   # M102: -869.7334000000001 -> -823.375
$ws.Cells.Item(102, 14).Value = -6825.3333   # N102: -7163.6 -> -6825.3333

# Row 107
$ws.Cells.Item(107, 8).Value = 72318.64   # H107: 77874.69500000001 -> 72318.64
$ws.Cells.Item(107, 9).Value = 143008.42   # I107: 166828.17 -> 143008.42
$ws.Cells.Item(107, 11).Value = 143008.42   # K107: 166828.17 -> 143008.42
$ws.Cells.Item(107, 13).Value = -141088.42   # M107: -164908.17 -> -141088.42

# Row 133
$ws.Cells.Item(133, 8).Value = 89994   # H133: 89994.5 -> 89994
$ws.Cells.Item(133, 10).Value = 89994   # J133: 89994.5 -> 89994
$ws.Cells.Item(133, 12).Value = 89994   # L133: 89994.5 -> 89994
$ws.Cells.Item(133, 14).Value = -100114   # N133: -100114.5 -> -100114

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Cells.Item(7, 8).Value = 4428.7144   # H7: 4087.5 -> 4428.7144
$ws.Cells.Item(7, 9).Value = 2200.6   # I7: 2117 -> 2200.6
$ws.Cells.Item(7, 11).Value = 2200.6   # K7: 2117 -> 2200.6
$ws.Cells.Item(7, 13).Value = -2088.6   # M7: -2005 -> -2088.6

# Row 40
$ws.Cells.Item(40, 8).Value = 1825   # H40: 1990.909 -> 1825
$ws.Cells.Item(40, 9).Value = 1036.3636   # I40: 1140 -> 1036.3636
$ws.Cells.Item(40, 11).Value = 1036.3636   # K40: 1140 -> 1036.3636
$ws.Cells.Item(40, 13).Value = -900.3635999999999   # M40: -1004 -> -900.3635999999999

# Row 46
$ws.Cells.Item(46, 8).Value = 5334.8335   # H46: 5337.528 -> 5334.8335
$ws.Cells.Item(46, 9).Value = 37595   # I46: 37596.332 -> 37595
$ws.Cells.Item(46, 10).Value = 2402.0908   # J46: 2404.9092 -> 2402.0908
$ws.Cells.Item(46, 11).Value = 37595   # K46: 37596.332 -> 37595
$ws.Cells.Item(46, 12).Value = 2402.0908   # L46: 2404.9092 -> 2402.0908
$ws.Cells.Item(46, 13).Value = -37407   # M46: -37408.332 -> -37407
$ws.Cells.Item(46, 14).Value = -2778.0908   # N46: -2780.9092 -> -2778.0908

# Row 63
$ws.Cells.Item(63, 8).Value = 0   # H63: 69999 -> 0
$ws.Cells.Item(63, 10).Value = 0   # J63: 69999 -> 0
$ws.Cells.Item(63, 12).Value = 0   # L63: 69999 -> 0
$ws.Cells.Item(63, 14).Value = ""   # N63: -71497 -> None

# Row 66
$ws.Cells.Item(66, 8).Value = 0   # H66: 69999 -> 0
$ws.Cells.Item(66, 10).Value = 0   # J66: 69999 -> 0
$ws.Cells.Item(66, 12).Value = 0   # L66: 209997 -> 0
$ws.Cells.Item(66, 14).Value = ""   # N66: -217485 -> None

# Row 102
$ws.Cells.Item(102, 8).Value = 69999   # H102: 0 -> 69999
$ws.Cells.Item(102, 10).Value = 69999   # J102: 0 -> 69999
$ws.Cells.Item(102, 12).Value = 69999   # L102: 0 -> 69999
$ws.Cells.Item(102, 14).Value = -76489   # N102: None -> -76489

# Row 126
$ws.Cells.Item(126, 8).Value = 4428.7144   # H126: 4087.5 -> 4428.7144
$ws.Cells.Item(126, 9).Value = 2200.6   # I126: 2117 -> 2200.6
$ws.Cells.Item(126, 11).Value = 6601.799999999999   # K126: 6351 -> 6601.799999999999
$ws.Cells.Item(126, 13).Value = -4131.799999999999   # M126: -3881 -> -4131.799999999999

# Row 132
$ws.Cells.Item(132, 8).Value = 96260.08   # H132: 89562.92999999999 -> 96260.08
$ws.Cells.Item(132, 9).Value = 103031.75   # I132: 95298.53999999999 -> 103031.75
$ws.Cells.Item(132, 11).Value = 309095.25   # K132: 285895.62 -> 309095.25
$ws.Cells.Item(132, 13).Value = -306565.25   # M132: -283365.62 -> -306565.25

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 43
$ws.Cells.Item(43, 8).Value = 54999   # H43: 0 -> 54999
$ws.Cells.Item(43, 10).Value = 54999   # J43: 0 -> 54999
$ws.Cells.Item(43, 12).Value = 54999   # L43: 0 -> 54999
$ws.Cells.Item(43, 14).Value = -55297   # N43: None -> -55297

# Row 70
$ws.Cells.Item(70, 8).Value = 0   # H70: 49999 -> 0
$ws.Cells.Item(70, 10).Value = 0   # J70: 49999 -> 0
$ws.Cells.Item(70, 12).Value = 0   # L70: 49999 -> 0
$ws.Cells.Item(70, 14).Value = ""   # N70: -50629 -> None

# Row 73
$ws.Cells.Item(73, 8).Value = 0   # H73: 49999 -> 0
$ws.Cells.Item(73, 10).Value = 0   # J73: 49999 -> 0
$ws.Cells.Item(73, 12).Value = 0   # L73: 49999 -> 0
$ws.Cells.Item(73, 14).Value = ""   # N73: -52183 -> None

# Row 81
$ws.Cells.Item(81, 8).Value = 1899.2858   # H81: 1892.1538 -> 1899.2858
$ws.Cells.Item(81, 10).Value = 2942.5   # J81: 3132.6 -> 2942.5
$ws.Cells.Item(81, 12).Value = 5885   # L81: 6265.2 -> 5885
$ws.Cells.Item(81, 14).Value = -8007   # N81: -8387.200000000001 -> -8007

# Row 84
$ws.Cells.Item(84, 8).Value = 1899.2858   # H84: 1892.1538 -> 1899.2858
$ws.Cells.Item(84, 10).Value = 2942.5   # J84: 3132.6 -> 2942.5
$ws.Cells.Item(84, 12).Value = 29425   # L84: 31326 -> 29425
$ws.Cells.Item(84, 14).Value = -40033   # N84: -41934 -> -40033

# Row 125
$ws.Cells.Item(125, 8).Value = 54632.668   # H125: 55501.43 -> 54632.668
$ws.Cells.Item(125, 10).Value = 54632.668   # J125: 55501.43 -> 54632.668
$ws.Cells.Item(125, 12).Value = 54632.668   # L125: 55501.43 -> 54632.668
$ws.Cells.Item(125, 14).Value = -64472.668   # N125: -65341.43 -> -64472.668

# Row 132
$ws.Cells.Item(132, 8).Value = 59277.332   # H132: 62617.234 -> 59277.332
$ws.Cells.Item(132, 9).Value = 59277.332   # I132: 62617.234 -> 59277.332
$ws.Cells.Item(132, 11).Value = 177831.996   # K132: 187851.702 -> 177831.996
$ws.Cells.Item(132, 13).Value = -175301.996   # M132: -185321.702 -> -175301.996

# Row 136
$ws.Cells.Item(136, 8).Value = 2109.2942   # H136: 2222.3125 -> 2109.2942
$ws.Cells.Item(136, 9).Value = 1709.875   # I136: 1803.8 -> 1709.875
$ws.Cells.Item(136, 11).Value = 5129.625   # K136: 5411.4 -> 5129.625
$ws.Cells.Item(136, 13).Value = -2579.625   # M136: -2861.4 -> -2579.625
